# Publish IG 1.0.1
# - rename "Include from unknown" sheet to "Include #0"
# - bump Version to 1.0.1
# - update Contact display value
# - add a new "Jurisdiction" metadata row (with an empty value) right after "Contact"

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsInclude = $wb.Worksheets.Item(2)

# 1) Rename the "Include from unknown" sheet to "Include #0"
$wsInclude.Name = "Include #0"

# 2) Version 1.0.0 -> 1.0.1
$wsMeta.Range("B3").Value = "1.0.1"

# 3) Contact value -> MedCom (http://www.medcom.dk)
$wsMeta.Range("B10").Value = "MedCom (http://www.medcom.dk)"

# 4) Insert a new "Jurisdiction" row right after the "Contact" row (row 10),
#    pushing Description/Purpose/Copyright/Immutable down by one row.
$wsMeta.Rows.Item(11).Insert()

# Match the formatting of the other data rows (border/wrap/style) by copying
# the style from the row right below (Description, which kept the regular
# data-row style) onto the freshly inserted row.
$wsMeta.Range("A12:B12").Copy()
$wsMeta.Range("A11:B11").PasteSpecial(-4122)

$wsMeta.Range("A11").Value = "Jurisdiction"
$wsMeta.Range("B11").Value = ""
